$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title
#    (first) paragraph. It must match the document's existing pattern of
#    body paragraphs: a leading empty run, then a bold "Meta description"
#    run, then a plain run with the rest of the sentence.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)

$metaXml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Experience an exciting and easy to play slot game with Dragon Gate Trial. Play for free with a high level of diversification in reels and features.</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$insertionPoint = $d.Range($p1.Range.End, $p1.Range.End)
[void]$insertionPoint.InsertXML($metaXml)

# InsertXML leaves behind an extra trailing empty paragraph (needed to make
# the break "stick") - remove it, which restores the following paragraph
# ("Get Lucky With Dragon Gate Trials Features") to its original state.
$extra = $d.Paragraphs.Item(3)
$extra.Range.Delete()

# ---------------------------------------------------------------------------
# 2. Near the end of the document: drop the duplicated bold title
#    paragraph, and rewrite the italic paragraph's text with the new
#    image-generation prompt (keeping its existing run/format structure).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldDupe = $d.Paragraphs.Item($count - 1)
$boldDupe.Range.Delete()

$count = $d.Paragraphs.Count
$italicPara = $d.Paragraphs.Item($count)
$newRange = $d.Range($italicPara.Range.Start, $italicPara.Range.End - 1)
$newRange.Text = "Please create a feature image for Dragon Gate Trial that fits the following criteria: - The image should be in cartoon style - The image should feature a happy Maya warrior with glasses. Note: The image should not include any references to the game's actual theme of Chinese New Year and dragons. The Maya warrior should be the main focus of the image."
